$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.9758980442694933, 1.1892427924274336, 0.314572324456495, 3.0340322574224636)"
$ws.Range("C2").Value = "NCT(2.5900620673138803, 1.7152424906894583, -0.012907075629254482, 4.587860877770693)"
$ws.Range("D2").Value = "NIG(2.8752316438471013, 1.9700193438479292, -0.566268220921977, 5.756398371846208)"
$ws.Range("E2").Value = "NIG(1.3008730560350323, 0.9650857246727593, 4.200028673643869, 5.235245364407481)"
